# Apply weekly price/date update for Hortaliza, Mapocho Venta Directa de Santiago - Alcachofa
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 11 data)
$ws.Range("D2").Value = 44432
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("O2").Value = "Provincia del Elquí"
$ws.Range("P2").Value = 467

# Row 3 (was row 5 data)
$ws.Range("D3").Value = 44474
$ws.Range("J3").Value = 45
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 333

# Row 4 (was row 2 data)
$ws.Range("D4").Value = 44453
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 12000
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 400

# Row 5 (was row 7 data)
$ws.Range("D5").Value = 44460
$ws.Range("J5").Value = 45
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 13000
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 433

# Row 6 (was row 12 data)
$ws.Range("D6").Value = 44435
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 14000
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 467

# Row 7 (was row 13 data)
$ws.Range("D7").Value = 44435
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 14000
$ws.Range("O7").Value = "Provincia del Elquí"
$ws.Range("P7").Value = 467

# Row 8 (was row 4 data)
$ws.Range("D8").Value = 44418
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 500

# Row 9 (was row 3 data)
$ws.Range("D9").Value = 44376
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 18000
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 600

# Row 10 (was row 8 data)
$ws.Range("D10").Value = 44449
$ws.Range("J10").Value = 45
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 12000
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 400

# Row 11 (was row 6 data)
$ws.Range("D11").Value = 44446
$ws.Range("J11").Value = 25
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 14000
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 467

# Row 12 (was row 9 data)
$ws.Range("D12").Value = 44425
$ws.Range("J12").Value = 35
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 14000
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 467

# Row 13 (was row 10 data)
$ws.Range("D13").Value = 44421
$ws.Range("J13").Value = 25
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15400
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 513

